$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7892.251036985298
$ws.Range("C2").Value = 1938.454135688329
$ws.Range("D2").Value = 13973.80148380539
